$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.036.99'
$ws.Range('E2').Value = '  -2.12%  '

$ws.Range('D3').Value = '2.491.77'
$ws.Range('E3').Value = '  -3.69%  '

$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.86'
$ws.Range('E5').Value = '  -0.38%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.57'
$ws.Range('E6').Value = '  -1.80%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.581'
$ws.Range('E7').Value = '  +0.74%  '

$ws.Range('E8').Value = '  +0.20%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.531'
$ws.Range('E9').Value = '  -3.57%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.19'
$ws.Range('E10').Value = '  -1.97%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0802'
$ws.Range('E11').Value = '  -1.04%  '

$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.59'
$ws.Range('E12').Value = '  -2.12%  '

$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.112'
$ws.Range('E13').Value = '  -2.01%  '

$ws.Range('D14').Value = '2.880.44'
$ws.Range('E14').Value = '  -3.36%  '

$ws.Range('D15').Value = '2.485.21'
$ws.Range('E15').Value = '  -4.30%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.84'
$ws.Range('E16').Value = '  +2.85%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.852'
$ws.Range('E17').Value = '  -4.32%  '

$ws.Range('D18').Value = '42.124.09'
$ws.Range('E18').Value = '  -2.08%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.78'
$ws.Range('E19').Value = '  -1.49%  '

$ws.Range('D20').Value = '0.0₃0964'
$ws.Range('E20').Value = '  -2.92%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.37'
$ws.Range('E21').Value = '  -4.73%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '70.67'
$ws.Range('E22').Value = '  -1.99%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '247.72'
$ws.Range('E23').Value = '  -3.02%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.88'
$ws.Range('E24').Value = '  -2.99%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.00'
$ws.Range('E25').Value = '  -6.63%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.27'
$ws.Range('E26').Value = '  -8.76%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.14%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.31'
$ws.Range('E28').Value = '  +8.96%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.09'
$ws.Range('E29').Value = '  -1.51%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '37.14'
$ws.Range('E30').Value = '  -4.65%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.87'
$ws.Range('E31').Value = '  -3.49%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '154.27'
$ws.Range('E32').Value = '  -0.92%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.28'
$ws.Range('E33').Value = '  -3.13%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0778'
$ws.Range('E34').Value = '  -4.53%  '

$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.61'
$ws.Range('E35').Value = '  -5.42%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.04'
$ws.Range('E36').Value = '  -6.88%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.34'
$ws.Range('E37').Value = '  -0.58%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.114'
$ws.Range('E38').Value = '  -0.61%  '

$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '23.93'
$ws.Range('E39').Value = '  +1.75%  '

$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.118'
$ws.Range('E40').Value = '  -1.52%  '

$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.81'
$ws.Range('E41').Value = '  -2.94%  '

$ws.Range('B42').Value = 'NEARProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.34'
$ws.Range('E42').Value = '  -2.74%  '

$ws.Range('E43').Value = '  +0.08%  '

$ws.Range('D44').Value = '2.039.29'
$ws.Range('E44').Value = '  -1.47%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0296'
$ws.Range('E45').Value = '  -5.04%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.94'
$ws.Range('E46').Value = '  -6.38%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '83.09'
$ws.Range('E47').Value = '  -2.90%  '

$ws.Range('E48').Value = '  -4.26%  '

$ws.Range('D49').Value = '2.750.15'
$ws.Range('E49').Value = '  -2.80%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.187'
$ws.Range('E50').Value = '  -2.65%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '100.41'
$ws.Range('E51').Value = '  -5.69%  '
